# Spring2015.xlsx - "Changed table generation procedures"
#
# 1. Make the "TableMappings" sheet the active tab/sheet (it was "Companies").
# 2. Move the selection/scroll position on "TableMappings" to F28 (top-left
#    near A22); "Companies" keeps its own topLeftCell/selection but is no
#    longer the tabSelected sheet.
# 3. Column A ("Location") on "TableMappings", rows 35-73, is renumbered to
#    continue the row-1 sequence (fixing a duplicate/skip in the old data).

$wb = $excel.ActiveWorkbook
$companies = $wb.Worksheets.Item("Companies")
$tableMappings = $wb.Worksheets.Item("TableMappings")

# --- Fix up the Location numbering on TableMappings (rows 35-73) ---
for ($r = 35; $r -le 73; $r++) {
    $tableMappings.Cells.Item($r, 1).Value = $r - 1
}

# --- Move the selection on TableMappings to its new spot ---
$tableMappings.Range("F28").Select()

# --- TableMappings becomes the active/selected sheet & tab ---
$tableMappings.Activate()
